$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '61.472.43'
Set-TextValue 'D3' '2.958.46'
Set-TextValue 'E3' '  -6.33%  '
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '542.51'
Set-TextValue 'E5' '  -4.94%  '
Set-TextValue 'D6' '151.75'
Set-TextValue 'E6' '  -7.38%  '
Set-TextValue 'E7' '  -0.13%  '
Set-TextValue 'E8' '  -1.34%  '
Set-TextValue 'D9' '2.967.01'
Set-TextValue 'E9' '  -6.09%  '
Set-TextValue 'E10' '  -3.31%  '
Set-TextValue 'E11' '  -7.28%  '
Set-TextValue 'E12' '  -3.59%  '
Set-TextValue 'D13' '3.475.66'
Set-TextValue 'E13' '  -6.31%  '
Set-TextValue 'E14' '  -2.83%  '
Set-TextValue 'D15' '61.554.01'
Set-TextValue 'E15' '  -4.45%  '
Set-TextValue 'E16' '  -5.82%  '
Set-TextValue 'B17' 'WrappedEther'
Set-TextValue 'C17' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D17' '2.960.38'
Set-TextValue 'E17' '  -6.43%  '
Set-TextValue 'B18' 'ShibaInu'
Set-TextValue 'C18' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D18' '0.0000147'
Set-TextValue 'E18' '  -5.01%  '
Set-TextValue 'D19' '5.17'
Set-TextValue 'E19' '  -1.22%  '
Set-TextValue 'D20' '381.91'
Set-TextValue 'E20' '  -5.68%  '
Set-TextValue 'D21' '11.99'
Set-TextValue 'E21' '  -5.38%  '
Set-TextValue 'D22' '6.68'
Set-TextValue 'E22' '  -6.07%  '
Set-TextValue 'E23' '  +0.04%  '
Set-TextValue 'D24' '65.28'
Set-TextValue 'E24' '  -4.75%  '
Set-TextValue 'E25' '  -2.90%  '
Set-TextValue 'D26' '3.083.52'
Set-TextValue 'E26' '  -6.63%  '
Set-TextValue 'D27' '0.186'
Set-TextValue 'E27' '  -3.95%  '
Set-TextValue 'D28' '0.997'
Set-TextValue 'E28' '  -0.20%  '
Set-TextValue 'D29' '0.0₃0933'
Set-TextValue 'E29' '  -8.36%  '
Set-TextValue 'E30' '  -5.49%  '
Set-TextValue 'E32' '  -5.06%  '
Set-TextValue 'E33' '  -3.50%  '
Set-TextValue 'D34' '159.33'
Set-TextValue 'E34' '  +1.61%  '
Set-TextValue 'E35' '  -3.59%  '
Set-TextValue 'E36' '  -5.26%  '
Set-TextValue 'E37' '  -4.68%  '
Set-TextValue 'E38' '  -4.83%  '
Set-TextValue 'E39' '  -7.14%  '
Set-TextValue 'E40' '  -3.50%  '
Set-TextValue 'D41' '2.412.25'
Set-TextValue 'E41' '  -9.67%  '
Set-TextValue 'D42' '37.22'
Set-TextValue 'E42' '  -3.36%  '
Set-TextValue 'D43' '22.15'
Set-TextValue 'E43' '  -6.99%  '
Set-TextValue 'D44' '0.663'
Set-TextValue 'E44' '  -4.60%  '
Set-TextValue 'D45' '0.0595'
Set-TextValue 'E45' '  -3.39%  '
Set-TextValue 'D46' '0.998'
Set-TextValue 'E46' '  -0.23%  '
Set-TextValue 'E47' '  -3.84%  '
Set-TextValue 'E48' '  -8.95%  '
Set-TextValue 'E49' '  -2.39%  '
Set-TextValue 'B50' 'InjectiveProtocol'
Set-TextValue 'C50' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D50' '19.78'
Set-TextValue 'E50' '  -6.97%  '
Set-TextValue 'B51' 'Bittensor'
Set-TextValue 'C51' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D51' '268.17'
Set-TextValue 'E51' '  -6.93%  '
